# Decrement column E (剩余) values by 1 for every data row (2-99),
# except row 36, which stays unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($r, 5)  # Column E = 5
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
